$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-09 21:18:54"
$ws.Range("I2").Value = "2.2 mm"
$ws.Range("E3").Value = "2026-02-09 21:18:56"
$ws.Range("G3").Value = "169 cm"
$ws.Range("I3").Value = "2.7 mm"
$ws.Range("E4").Value = "2026-02-09 21:18:59"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "85%"
$ws.Range("O4").Value = "8.2 °C"
$ws.Range("E5").Value = "2026-02-09 21:19:02"
$ws.Range("E6").Value = "2026-02-09 21:19:04"
$ws.Range("E7").Value = "2026-02-09 21:19:07"
$ws.Range("E8").Value = "2026-02-09 21:19:10"
$ws.Range("E9").Value = "2026-02-09 21:19:12"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "82%"
$ws.Range("O9").Value = "8.5 °C"
$ws.Range("E10").Value = "2026-02-09 21:19:15"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "84%"
$ws.Range("E11").Value = "2026-02-09 21:19:18"
$ws.Range("E12").Value = "2026-02-09 21:19:20"
$ws.Range("E13").Value = "2026-02-09 21:19:22"
$ws.Range("E14").Value = "2026-02-09 21:19:25"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "77%"
$ws.Range("E15").Value = "2026-02-09 21:19:28"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "81%"
$ws.Range("E16").Value = "2026-02-09 21:19:30"
$ws.Range("I16").Value = "2.2 mm"
$ws.Range("E17").Value = "2026-02-09 21:19:33"
$ws.Range("E18").Value = "2026-02-09 21:19:36"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "82%"
$ws.Range("E19").Value = "2026-02-09 21:19:38"
$ws.Range("E20").Value = "2026-02-09 21:19:41"
$ws.Range("I20").Value = "0.4 mm"
$ws.Range("E21").Value = "2026-02-09 21:19:43"
$ws.Range("I21").Value = "0.1 mm"
$ws.Range("E22").Value = "2026-02-09 21:19:46"
$ws.Range("E23").Value = "2026-02-09 21:19:49"
$ws.Range("I23").Value = "1.5 mm"
$ws.Range("E24").Value = "2026-02-09 21:19:51"
$ws.Range("E25").Value = "2026-02-09 21:19:54"
$ws.Range("I25").Value = "0.3 mm"
$ws.Range("E26").Value = "2026-02-09 21:19:57"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "84%"
$ws.Range("J26").Value = "1006.7 hPa"
$ws.Range("O26").Value = "2.6 °C"
$ws.Range("E27").Value = "2026-02-09 21:19:59"
$ws.Range("I27").Value = "0.5 mm"
$ws.Range("O27").Value = "-2.3 °C"
$ws.Range("E28").Value = "2026-02-09 21:20:02"
$ws.Range("O28").Value = "7.3 °C"
$ws.Range("E29").Value = "2026-02-09 21:20:05"
$ws.Range("O29").Value = "8.5 °C"
$ws.Range("E30").Value = "2026-02-09 21:20:07"
$ws.Range("E31").Value = "2026-02-09 21:20:10"
$ws.Range("J31").Value = "1006.3 hPa"
$ws.Range("E32").Value = "2026-02-09 21:20:13"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "81%"
$ws.Range("I32").Value = "1.0 mm"
$ws.Range("E33").Value = "2026-02-09 21:20:16"
$ws.Range("E34").Value = "2026-02-09 21:20:18"
$ws.Range("E35").Value = "2026-02-09 21:20:21"
$ws.Range("E36").Value = "2026-02-09 21:20:23"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "79%"
$ws.Range("O36").Value = "9.7 °C"
$ws.Range("E37").Value = "2026-02-09 21:20:26"
$ws.Range("E38").Value = "2026-02-09 21:20:29"
$ws.Range("E39").Value = "2026-02-09 21:20:31"
$ws.Range("L39").Value = "52.6 km/h - 299º 20:54 TU"
$ws.Range("O39").Value = "-3.3 °C"
$ws.Range("E40").Value = "2026-02-09 21:20:34"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "82%"
$ws.Range("I40").Value = "0.4 mm"
$ws.Range("E41").Value = "2026-02-09 21:20:36"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "57%"
$ws.Range("E42").Value = "2026-02-09 21:20:39"
$ws.Range("E43").Value = "2026-02-09 21:20:42"
$ws.Range("E44").Value = "2026-02-09 21:20:44"
$ws.Range("I44").Value = "1.2 mm"
$ws.Range("E45").Value = "2026-02-09 21:20:47"
$ws.Range("I45").Value = "1.0 mm"
$ws.Range("E46").Value = "2026-02-09 21:20:49"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "74%"
$ws.Range("I46").Value = "0.9 mm"
